$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.711363077163696
$ws.Range("B1").Value = 1.909371376037598
$ws.Range("C1").Value = 1.9500732421875
$ws.Range("D1").Value = 2.253790140151978
$ws.Range("E1").Value = 3.001873016357422
